$wb = $excel.ActiveWorkbook

# --- Sheet1: update a few numeric-looking text values ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A1").Value = "9,19"
$ws1.Range("A2").Value = "4,3748"
$ws1.Range("A43").Value = "1672,69"
$ws1.Range("A213").Value = "26937,03"

# --- data sheet: update headers and clear stale "Sheet1"/cell-ref rows ---
$wsData = $wb.Worksheets.Item("data")
$wsData.Range("A1").Value = "avalanche"
$wsData.Range("B1").Value = "-"
$wsData.Range("C1").Value = "-"

$wsData.Range("B2").Value = ""
$wsData.Range("C2").Value = ""

$wsData.Range("A3").Value = "A1"
$wsData.Range("B3").Value = ""
$wsData.Range("C3").Value = ""
